$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-9 in the "Fitness" column (C) lose their value (1.0) and their
# number-format style, leaving bare empty cells (<c r="C3"/> etc.),
# while row 2's C2 is left untouched.
$target = $ws.Range("C3:C9")
$target.ClearContents()
$target.ClearFormats()
